$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 57; existing rows 57..126 shift down to 58..127
$ws.Rows(57).Insert()

# Populate the newly inserted row 57 with the new record's data
$ws.Cells.Item(57, 1).Value = 7
$ws.Cells.Item(57, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(57, 3).Value = "Ñuble"
$ws.Cells.Item(57, 4).Value = 45079
$ws.Cells.Item(57, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(57, 5).Value = 16
$ws.Cells.Item(57, 6).Value = 100112030
$ws.Cells.Item(57, 7).Value = "Poroto granado"
$ws.Cells.Item(57, 8).Value = "Sin especificar"
$ws.Cells.Item(57, 9).Value = "Primera"
$ws.Cells.Item(57, 10).Value = 50
$ws.Cells.Item(57, 11).Value = 25000
$ws.Cells.Item(57, 12).Value = 25000
$ws.Cells.Item(57, 13).Value = 25000
$ws.Cells.Item(57, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(57, 15).Value = "Provincia de Diguillín"
$ws.Cells.Item(57, 16).Value = 1000
$ws.Cells.Item(57, 17).Value = 25
$ws.Cells.Item(57, 18).Value = "Hortaliza"
